$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1293.3334  # was 1205.5714
$ws.Range("I70").Value = 2000  # was 1259.6666
$ws.Range("J70").Value = 1152  # was 1165
$ws.Range("K70").Value = 6000  # was 3778.9998
$ws.Range("L70").Value = 3456  # was 3495
$ws.Range("M70").Value = -5730  # was -3508.9998
$ws.Range("N70").Value = -3996  # was -4035
$ws.Range("H73").Value = 1293.3334  # was 1205.5714
$ws.Range("I73").Value = 2000  # was 1259.6666
$ws.Range("J73").Value = 1152  # was 1165
$ws.Range("K73").Value = 6000  # was 3778.9998
$ws.Range("L73").Value = 3456  # was 3495
$ws.Range("M73").Value = -5064  # was -2842.9998
$ws.Range("N73").Value = -5328  # was -5367
$ws.Range("H94").Value = 1999.5  # was 2000
$ws.Range("I94").Value = 1999.5  # was 2000
$ws.Range("K94").Value = 1999.5  # was 2000
$ws.Range("M94").Value = -1548.5  # was -1549
$ws.Range("H110").Value = 29963.334  # was 29890
$ws.Range("J110").Value = 29963.334  # was 29890
$ws.Range("L110").Value = 29963.334  # was 29890
$ws.Range("N110").Value = -38143.334  # was -38070
$ws.Range("H132").Value = 2451.7026  # was 2749.697
$ws.Range("I132").Value = 2451.7026  # was 2749.697
$ws.Range("K132").Value = 7355.1078  # was 8249.091
$ws.Range("M132").Value = -4825.1078  # was -5719.091
$ws.Range("H138").Value = 135838.88  # was 147348.92
$ws.Range("I138").Value = 3093.75  # was 4570
$ws.Range("J138").Value = 151455.95  # was 158331.9
$ws.Range("K138").Value = 9281.25  # was 13710
$ws.Range("L138").Value = 454367.85  # was 474995.7
$ws.Range("M138").Value = -4141.25  # was -8570
$ws.Range("N138").Value = -464647.85  # was -485275.7

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 4000  # was 1309.5
$ws.Range("I30").Value = 0  # was 609
$ws.Range("J30").Value = 4000  # was 2010
$ws.Range("K30").Value = 0  # was 609
$ws.Range("L30").Value = 4000  # was 2010
$ws.Range("M30").ClearContents()  # was -459
$ws.Range("N30").Value = -4300  # was -2310
$ws.Range("H32").Value = 9432.826999999999  # was 8020.89
$ws.Range("I32").Value = 6593.164  # was 5514.5713
$ws.Range("K32").Value = 6593.164  # was 5514.5713
$ws.Range("M32").Value = -6306.164  # was -5227.5713
$ws.Range("H74").Value = 21277798  # was 24391610
$ws.Range("I74").Value = 25641588  # was 30303646
$ws.Range("J74").Value = 4330  # was 4462.5
$ws.Range("K74").Value = 25641588  # was 30303646
$ws.Range("L74").Value = 4330  # was 4462.5
$ws.Range("M74").Value = -25640714  # was -30302772
$ws.Range("N74").Value = -6078  # was -6210.5
$ws.Range("H77").Value = 21277798  # was 24391610
$ws.Range("I77").Value = 25641588  # was 30303646
$ws.Range("J77").Value = 4330  # was 4462.5
$ws.Range("K77").Value = 128207940  # was 151518230
$ws.Range("L77").Value = 21650  # was 22312.5
$ws.Range("M77").Value = -128203572  # was -151513862
$ws.Range("N77").Value = -30386  # was -31048.5
$ws.Range("H102").Value = 1413.25  # was 1388.8096
$ws.Range("I102").Value = 1309.7059  # was 1286.9445
$ws.Range("K102").Value = 1309.7059  # was 1286.9445
$ws.Range("M102").Value = 312.2941000000001  # was 335.0554999999999
$ws.Range("H132").Value = 12265.9795  # was 13436.796
$ws.Range("I132").Value = 1721.7805  # was 1810.3422
$ws.Range("J132").Value = 66305  # was 87071
$ws.Range("K132").Value = 5165.3415  # was 5431.0266
$ws.Range("L132").Value = 198915  # was 261213
$ws.Range("M132").Value = -2635.3415  # was -2901.0266
$ws.Range("N132").Value = -203975  # was -266273

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 862.5172  # was 878.3570999999999
$ws.Range("J80").Value = 617.6111  # was 629.2941
$ws.Range("L80").Value = 617.6111  # was 629.2941
$ws.Range("N80").Value = -2613.6111  # was -2625.2941
$ws.Range("H83").Value = 862.5172  # was 878.3570999999999
$ws.Range("J83").Value = 617.6111  # was 629.2941
$ws.Range("L83").Value = 3088.0555  # was 3146.4705
$ws.Range("N83").Value = -13072.0555  # was -13130.4705
$ws.Range("H86").Value = 1920.7916  # was 1919.3334
$ws.Range("I86").Value = 1768.6875  # was 1776.4117
$ws.Range("J86").Value = 2225  # was 2526.75
$ws.Range("K86").Value = 1768.6875  # was 1776.4117
$ws.Range("L86").Value = 2225  # was 2526.75
$ws.Range("M86").Value = -645.6875  # was -653.4117000000001
$ws.Range("N86").Value = -4471  # was -4772.75
$ws.Range("H89").Value = 1920.7916  # was 1919.3334
$ws.Range("I89").Value = 1768.6875  # was 1776.4117
$ws.Range("J89").Value = 2225  # was 2526.75
$ws.Range("K89").Value = 8843.4375  # was 8882.058500000001
$ws.Range("L89").Value = 11125  # was 12633.75
$ws.Range("M89").Value = -3227.4375  # was -3266.058500000001
$ws.Range("N89").Value = -22357  # was -23865.75
$ws.Range("H105").Value = 1113131.1  # was 1221604.9
$ws.Range("I105").Value = 1638.125  # was 1679.2858
$ws.Range("J105").Value = 1726368.6  # was 1854158.9
$ws.Range("K105").Value = 1638.125  # was 1679.2858
$ws.Range("L105").Value = 1726368.6  # was 1854158.9
$ws.Range("M105").Value = 108.875  # was 67.71419999999989
$ws.Range("N105").Value = -1729862.6  # was -1857652.9

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 37184.5  # was 35016.332
$ws.Range("J18").Value = 37184.5  # was 35016.332
$ws.Range("L18").Value = 37184.5  # was 35016.332
$ws.Range("N18").Value = -37644.5  # was -35476.332
$ws.Range("H48").Value = 9800  # was 8400
$ws.Range("J48").Value = 9800  # was 8400
$ws.Range("L48").Value = 9800  # was 8400
$ws.Range("N48").Value = -10752  # was -9352
$ws.Range("H58").Value = 26585.6  # was 31159.53
$ws.Range("I58").Value = 1527.5  # was 1762.2727
$ws.Range("K58").Value = 1527.5  # was 1762.2727
$ws.Range("M58").Value = -1324.5  # was -1559.2727
$ws.Range("H60").Value = 19886.666  # was 17870.428
$ws.Range("I60").Value = 20000  # was 10046.5
$ws.Range("J60").Value = 19872.5  # was 21000
$ws.Range("K60").Value = 20000  # was 10046.5
$ws.Range("L60").Value = 19872.5  # was 21000
$ws.Range("M60").Value = -19489  # was -9535.5
$ws.Range("N60").Value = -20894.5  # was -22022
$ws.Range("H68").Value = 54625  # was 56556.25
$ws.Range("J68").Value = 54625  # was 56556.25
$ws.Range("L68").Value = 54625  # was 56556.25
$ws.Range("N68").Value = -56123  # was -58054.25
$ws.Range("H71").Value = 54625  # was 56556.25
$ws.Range("J71").Value = 54625  # was 56556.25
$ws.Range("L71").Value = 163875  # was 169668.75
$ws.Range("N71").Value = -171363  # was -177156.75
$ws.Range("H76").Value = 16669467  # was 16669833
$ws.Range("I76").Value = 16669467  # was 16669833
$ws.Range("K76").Value = 16669467  # was 16669833
$ws.Range("M76").Value = -16669152  # was -16669518
$ws.Range("H79").Value = 16669467  # was 16669833
$ws.Range("I79").Value = 16669467  # was 16669833
$ws.Range("K79").Value = 16669467  # was 16669833
$ws.Range("M79").Value = -16668375  # was -16668741
$ws.Range("H99").Value = 18521770  # was 18522040
$ws.Range("I99").Value = 2886.4666  # was 3046
$ws.Range("J99").Value = 41670376  # was 35718250
$ws.Range("K99").Value = 2886.4666  # was 3046
$ws.Range("L99").Value = 41670376  # was 35718250
$ws.Range("M99").Value = -1388.4666  # was -1548
$ws.Range("N99").Value = -41673372  # was -35721246
$ws.Range("H109").Value = 165021400  # was 165021580
$ws.Range("J109").Value = 165021400  # was 165021580
$ws.Range("L109").Value = 165021400  # was 165021580
$ws.Range("N109").Value = -165023480  # was -165023660
$ws.Range("H126").Value = 18521770  # was 18522040
$ws.Range("I126").Value = 2886.4666  # was 3046
$ws.Range("J126").Value = 41670376  # was 35718250
$ws.Range("K126").Value = 8659.399800000001  # was 9138
$ws.Range("L126").Value = 125011128  # was 107154750
$ws.Range("M126").Value = -6189.399800000001  # was -6668
$ws.Range("N126").Value = -125016068  # was -107159690
$ws.Range("H136").Value = 26585.6  # was 31159.53
$ws.Range("I136").Value = 1527.5  # was 1762.2727
$ws.Range("K136").Value = 4582.5  # was 5286.8181
$ws.Range("M136").Value = -2032.5  # was -2736.8181

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2450  # was 2433.3333
$ws.Range("J69").Value = 2500  # was 2450
$ws.Range("L69").Value = 7500  # was 7350
$ws.Range("N69").Value = -9122  # was -8972
$ws.Range("H72").Value = 2450  # was 2433.3333
$ws.Range("J72").Value = 2500  # was 2450
$ws.Range("L72").Value = 22500  # was 22050
$ws.Range("N72").Value = -30612  # was -30162
$ws.Range("H131").Value = 622.88776  # was 622.7732
$ws.Range("J131").Value = 749.9722  # was 751.60565
$ws.Range("L131").Value = 2249.9166  # was 2254.81695
$ws.Range("N131").Value = -12329.9166  # was -12334.81695
$ws.Range("H136").Value = 3536.348  # was 3413.52
$ws.Range("I136").Value = 1330  # was 1445.4546
$ws.Range("J136").Value = 4954.7144  # was 4959.857
$ws.Range("K136").Value = 3990  # was 4336.3638
$ws.Range("L136").Value = 14864.1432  # was 14879.571
$ws.Range("M136").Value = 1110  # was 763.6361999999999
$ws.Range("N136").Value = -25064.1432  # was -25079.571

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 24998.75  # was 0
$ws.Range("J108").Value = 24998.75  # was 0
$ws.Range("L108").Value = 24998.75  # was 0
$ws.Range("N108").Value = -32678.75  # was None

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3650.5  # was 3139.8
$ws.Range("I22").Value = 5451  # was 10001
$ws.Range("J22").Value = 1850  # was 1424.5
$ws.Range("K22").Value = 5451  # was 10001
$ws.Range("L22").Value = 1850  # was 1424.5
$ws.Range("M22").Value = -5156  # was -9706
$ws.Range("N22").Value = -2440  # was -2014.5
$ws.Range("H27").Value = 3650.5  # was 3139.8
$ws.Range("I27").Value = 5451  # was 10001
$ws.Range("J27").Value = 1850  # was 1424.5
$ws.Range("K27").Value = 5451  # was 10001
$ws.Range("L27").Value = 1850  # was 1424.5
$ws.Range("M27").Value = -5344  # was -9894
$ws.Range("N27").Value = -2064  # was -1638.5
$ws.Range("H68").Value = 2410.5557  # was 2053.3845
$ws.Range("I68").Value = 1566.6666  # was 1339.8
$ws.Range("J68").Value = 2832.5  # was 2499.375
$ws.Range("K68").Value = 1566.6666  # was 1339.8
$ws.Range("L68").Value = 2832.5  # was 2499.375
$ws.Range("M68").Value = -817.6666  # was -590.8
$ws.Range("N68").Value = -4330.5  # was -3997.375
$ws.Range("H71").Value = 2410.5557  # was 2053.3845
$ws.Range("I71").Value = 1566.6666  # was 1339.8
$ws.Range("J71").Value = 2832.5  # was 2499.375
$ws.Range("K71").Value = 7833.333000000001  # was 6699
$ws.Range("L71").Value = 14162.5  # was 12496.875
$ws.Range("M71").Value = -4089.333000000001  # was -2955
$ws.Range("N71").Value = -21650.5  # was -19984.875
$ws.Range("H82").Value = 1471.8125  # was 1425.2354
$ws.Range("I82").Value = 1474.5385  # was 1415.2667
$ws.Range("J82").Value = 1460  # was 1500
$ws.Range("K82").Value = 1474.5385  # was 1415.2667
$ws.Range("L82").Value = 1460  # was 1500
$ws.Range("M82").Value = -1113.5385  # was -1054.2667
$ws.Range("N82").Value = -2182  # was -2222
$ws.Range("H85").Value = 1471.8125  # was 1425.2354
$ws.Range("I85").Value = 1474.5385  # was 1415.2667
$ws.Range("J85").Value = 1460  # was 1500
$ws.Range("K85").Value = 1474.5385  # was 1415.2667
$ws.Range("L85").Value = 1460  # was 1500
$ws.Range("M85").Value = -226.5385000000001  # was -167.2666999999999
$ws.Range("N85").Value = -3956  # was -3996
$ws.Range("H93").Value = 2060.7  # was 2291.1765
$ws.Range("I93").Value = 1894.9412  # was 2139.2856
$ws.Range("K93").Value = 1894.9412  # was 2139.2856
$ws.Range("M93").Value = -646.9412  # was -891.2856000000002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3913.6365  # was 4366.6665
$ws.Range("I62").Value = 2920  # was 3500
$ws.Range("J62").Value = 4741.6665  # was 4800
$ws.Range("K62").Value = 2920  # was 3500
$ws.Range("L62").Value = 4741.6665  # was 4800
$ws.Range("M62").Value = -2296  # was -2876
$ws.Range("N62").Value = -5989.6665  # was -6048
$ws.Range("H65").Value = 3913.6365  # was 4366.6665
$ws.Range("I65").Value = 2920  # was 3500
$ws.Range("J65").Value = 4741.6665  # was 4800
$ws.Range("K65").Value = 14600  # was 17500
$ws.Range("L65").Value = 23708.3325  # was 24000
$ws.Range("M65").Value = -11480  # was -14380
$ws.Range("N65").Value = -29948.3325  # was -30240
